$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the consolidated dataset.
# Insert a new row at position 104 (pushes existing rows 104..183 down to 105..184)
# and fill it in with the new record's data. The descriptive columns
# (market/region/product taxonomy) are identical across every row in this sheet.
$ws.Rows(104).Insert()

$ws.Cells.Item(104, 1).Value = 4
$ws.Cells.Item(104, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(104, 3).Value = "Los Lagos"
$ws.Cells.Item(104, 4).Value = [DateTime]"2022-03-31"
$ws.Cells.Item(104, 5).Value = 10
$ws.Cells.Item(104, 6).Value = "Fruta"
$ws.Cells.Item(104, 7).Value = 100108
$ws.Cells.Item(104, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(104, 9).Value = 100108002
$ws.Cells.Item(104, 10).Value = "Mango"
$ws.Cells.Item(104, 11).Value = "Sin especificar"
$ws.Cells.Item(104, 12).Value = "Primera"
$ws.Cells.Item(104, 13).Value = 200
$ws.Cells.Item(104, 14).Value = 8000
$ws.Cells.Item(104, 15).Value = 8500
$ws.Cells.Item(104, 16).Value = 8250
$ws.Cells.Item(104, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(104, 18).Value = "Perú"
$ws.Cells.Item(104, 19).Value = 2062
$ws.Cells.Item(104, 20).Value = 4
